$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 <-> 17 swap (Chainlink/TRON reordered) with refreshed data ---
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D16")
$c.Value = "'0.127"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D17")
$c.Value = "'19.24"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.57%  "

# --- Price (D) / Volume(1h) (E) refreshes ---
$c = $ws.Range("D2")
$c.Value = "'68.958.25"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "

$c = $ws.Range("D3")
$c.Value = "'3.711.09"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  +0.06%  "

$c = $ws.Range("D5")
$c.Value = "'610.86"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.82%  "

$c = $ws.Range("D6")
$c.Value = "'188.96"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.63%  "

$c = $ws.Range("D7")
$c.Value = "'0.635"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("E8").Value = "  +0.29%  "

$c = $ws.Range("D9")
$c.Value = "'0.716"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.06%  "

$c = $ws.Range("D10")
$c.Value = "'0.159"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.88%  "

$c = $ws.Range("D11")
$c.Value = "'57.98"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.10%  "

$c = $ws.Range("D12")
$c.Value = "'0.0000288"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -4.73%  "

$c = $ws.Range("D13")
$c.Value = "'10.58"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.80%  "

$c = $ws.Range("D14")
$c.Value = "'4.312.28"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "

$c = $ws.Range("D15")
$c.Value = "'3.723.64"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("E18").Value = "  -1.44%  "

$c = $ws.Range("D19")
$c.Value = "'12.88"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.95%  "

$c = $ws.Range("D20")
$c.Value = "'68.790.46"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "

$c = $ws.Range("D21")
$c.Value = "'409.90"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.60%  "

$c = $ws.Range("D22")
$c.Value = "'4.59"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "

$c = $ws.Range("D23")
$c.Value = "'89.08"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "

$c = $ws.Range("D24")
$c.Value = "'3.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.63%  "

$c = $ws.Range("D25")
$c.Value = "'12.83"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "

$c = $ws.Range("D26")
$c.Value = "'10.85"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("E27").Value = "  +0.98%  "

$c = $ws.Range("D28")
$c.Value = "'3.78"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.28%  "

$c = $ws.Range("D29")
$c.Value = "'9.61"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.67%  "

$c = $ws.Range("D30")
$c.Value = "'32.98"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "

$c = $ws.Range("D31")
$c.Value = "'7.47"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -9.10%  "

$c = $ws.Range("D32")
$c.Value = "'12.68"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.05%  "

$c = $ws.Range("D33")
$c.Value = "'0.122"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.71%  "

$c = $ws.Range("D34")
$c.Value = "'45.85"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.33%  "

$c = $ws.Range("D35")
$c.Value = "'638.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.33%  "

$c = $ws.Range("D36")
$c.Value = "'65.65"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "

$c = $ws.Range("D37")
$c.Value = "'0.412"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.27%  "

$c = $ws.Range("D38")
$c.Value = "'0.0₃0815"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -13.32%  "

$c = $ws.Range("D39")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("E41").Value = "  +1.48%  "

$c = $ws.Range("D42")
$c.Value = "'3.03"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.86%  "

$c = $ws.Range("D43")
$c.Value = "'0.0444"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "

$c = $ws.Range("D44")
$c.Value = "'2.61"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "

$c = $ws.Range("D45")
$c.Value = "'0.139"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.41%  "

$c = $ws.Range("D46")
$c.Value = "'2.847.73"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.34%  "

$ws.Range("E47").Value = "  -0.37%  "

$c = $ws.Range("D48")
$c.Value = "'9.07"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.85%  "

$c = $ws.Range("D51")
$c.Value = "'2.57"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -22.73%  "

# --- Row 49 <-> 50 swap (Monero/ApeXProtocol reordered) with refreshed data ---
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D49")
$c.Value = "'3.11"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D50")
$c.Value = "'141.71"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.72%  "

